$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A ("carts" table, was "receipt_detail"/"receipt_id"/"product_id"/"quantity") ---
$ws.Range("A9").Value = "carts"

# A10 ("receipt_id", bold style) loses its bold styling and becomes "id"
$ws.Range("A10").Style = "Normal"
$ws.Range("A10").Value = "id"

$ws.Range("A11").Value = "quantity"
$ws.Range("A12").Value = "customer_id "

# New row: A13 "product_id " — copy formatting from an existing plain data cell (A12/A11 style)
$ws.Range("A11").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "product_id "

# --- Column E ("receipts" table header + fields) ---
$ws.Range("E5").Value = "date_of_birth"
$ws.Range("E7").Value = "password"

# E8 ("address") is no longer needed
$ws.Range("E8").Clear()

# E9 header becomes "receipts" and needs the same bold header styling as C9 ("admins")
$ws.Range("C9").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").Value = "receipts"

# New rows E10/E11
$ws.Range("E7").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Value = "id"

$ws.Range("E7").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E11").Value = "customer_id "

# E12 was the "receipts" header (bold style 3); it becomes a regular data
# row ("order_time") so it needs the plain data-cell styling back.
$ws.Range("E7").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("E12").Value = "order_time"
$ws.Range("E13").Value = "receiver_name"
$ws.Range("E14").Value = "receiver_phone"
$ws.Range("E15").Value = "receiver_address"
$ws.Range("E16").Value = "note"

# Trailing rows no longer used
$ws.Range("E17").Clear()
$ws.Range("E18").Clear()
$ws.Range("E19").Clear()

# --- View metadata: zoom 145 -> 130, selection moves to A16, scroll position resets ---
$excel.ActiveWindow.Zoom = 130
$ws.Range("A16").Select() | Out-Null
